$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.427132964134216
$ws.Range("B1").Value = 3.364690065383911
$ws.Range("C1").Value = 3.028459310531616
$ws.Range("D1").Value = 1.628336071968079
$ws.Range("E1").Value = 1.054091334342957
